$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "lookup" column D (and overflow column E for double matches) giving the
# raw source-list spelling that corresponds to each ISO "Code" in column C.
$ws.Range("D2").Value  = "x"
$ws.Range("D3").Value  = "CAPE VERDE"
$ws.Range("D4").Value  = "HOLY SEE (VATICAN CITY STATE)"
$ws.Range("D5").Value  = "Czech Republic"
$ws.Range("D6").Value  = "Cote D'ivoire"
$ws.Range("D7").Value  = "Iran, Islamic Republic of"
$ws.Range("D8").Value  = "KOREA, REPUBLIC OF"
$ws.Range("E8").Value  = "KOREA, DEMOCRATIC PEOPLE'S REPUBLIC OF"
$ws.Range("D9").Value  = "Libyan Arab Jamahiriya"
$ws.Range("D10").Value = "MICRONESIA, FEDERATED STATES OF"
$ws.Range("D11").Value = "x"
$ws.Range("D12").Value = "Saint Helena"
$ws.Range("D13").Value = "REUNION"
$ws.Range("D15").Value = "United States"
$ws.Range("D17").Value = "VIRGIN ISLANDS, BRITISH"
$ws.Range("E17").Value = "VIRGIN ISLANDS, U.S."

# D13 picked up a (cosmetic) explicit-font style during the original edit.
$ws.Range("D13").Font.Name = $ws.Range("D13").Font.Name

# Selection moved from D14 to D13 and the frozen/scrolled top-left cell reset.
$ws.Range("D13").Select()

# The printer resolution recorded in the page setup bumped on resave.
$ws.PageSetup.PrintQuality = 300
